# Thanh - weekly update
# Bump "Plan End" (K) and "Actual End" (M) dates for tasks 34-38 (rows 43-47)
# from 12/20/2016 to 12/21/2016, and move the active selection/viewport to
# reflect where the author was working (row ~33 area, cell M43).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($row in 43..47) {
    $ws.Cells.Item($row, 11).Value2 = 42725   # column K = Plan End
    $ws.Cells.Item($row, 13).Value2 = 42725   # column M = Actual End
}

# Reflect the author's on-screen context: scrolled down a row, and the
# cursor left on M43 after editing the last updated row.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 33
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("M43").Select()
